$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A: index, B: timestamp, C..L: counts). Row 1 header stays
# the same except column M ("RC6") is removed entirely below.
$data = @(
    @(0,  45393.33333333334, 1, 27, 2, 6, 1, 0, 105, 9, 31, 5),
    @(1,  45393.34027777778, 1, 32, 0, 1, 1, 3, 100, 5, 32, 2),
    @(2,  45393.34722222222, 1, 32, 0, 11, 1, 0, 92, 6, 31, 2),
    @(3,  45393.35416666666, 0, 15, 4, 11, 2, 2, 116, 4, 29, 1),
    @(4,  45393.36111111111, 0, 19, 3, 5, 4, 1, 114, 5, 25, 0),
    @(5,  45393.36805555555, 2, 27, 0, 10, 1, 0, 92, 1, 26, 3),
    @(6,  45393.66666666666, 0, 22, 0, 3, 0, 1, 105, 0, 17, 5),
    @(7,  45393.67361111111, 2, 21, 0, 2, 2, 1, 63, 0, 13, 3),
    @(8,  45393.68055555555, 0, 27, 0, 3, 0, 5, 71, 0, 9, 3),
    @(9,  45393.6875,        0, 39, 0, 2, 1, 1, 75, 0, 11, 2),
    @(10, 45393.69444444445, 0, 13, 0, 2, 3, 1, 58, 0, 16, 1),
    @(11, 45393.70138888889, 2, 35, 0, 5, 0, 2, 75, 0, 17, 5),
    @(12, 45393.83333333334, 0, 10, 0, 5, 0, 0, 43, 0, 17, 0),
    @(13, 45393.84027777778, 0, 11, 0, 2, 0, 1, 49, 0, 13, 0),
    @(14, 45393.84722222222, 0, 13, 0, 3, 0, 0, 60, 0, 12, 3),
    @(15, 45393.85416666666, 0, 22, 0, 5, 0, 0, 58, 0, 10, 0),
    @(16, 45393.86111111111, 0, 15, 0, 8, 2, 0, 50, 0, 7, 0),
    @(17, 45393.86805555555, 1, 21, 0, 4, 0, 0, 52, 0, 9, 0)
)

# Remove the RC6 column (M) entirely - collapses the used range back to
# column L, matching the new dimension A1:L19.
$ws.Columns.Item(13).Delete()

# The new data has 18 rows (was 12), so extend the formatted range (style
# carried by columns A/B for rows 2-13) down through row 19 before writing
# values, so the new rows pick up the same cell styles (index, date fmt).
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B19").PasteSpecial(-4122)

$row = 2
foreach ($r in $data) {
    for ($c = 0; $c -lt $r.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $r[$c]
    }
    $row++
}
